$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old sub-header row (row 2: "(m3/s)" / "(MW)" / "(GWh)" labels).
# Rows 3-14 shift up to become rows 2-13.
$ws.Rows.Item(2).Delete()

# Rebuild row 1 as a full header row across A1:K1.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# A1:E1 use the plain default style (no special font/format).
$ws.Range("A1:E1").ClearFormats()

# F1:K1 keep the small-font header style used elsewhere in the sheet.
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").Font.Name = "Arial"

# Match the selection left behind by the edit (A2:K2 highlighted).
$ws.Range("A2:K2").Select()
